$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grading update for row 15 (Dylan Zeledon): Homework 3, Quiz 1, Homework 4
$ws.Range("F15").Formula = "=28/30"
$ws.Range("G15").Formula = "=50/50"
$ws.Range("H15").Formula = "=20/20"

# Move the active selection to reflect where the editor left off
[void]$ws.Range("G16").Select()
